$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set Price column (D) cells as text to avoid Excel auto-number conversion
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.680.78'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.873.81'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '332.21'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4718'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3945'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '47.79'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08028'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.026'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.00'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.866.59'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.964'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.132'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.007'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.08'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06667'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.25'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.697.08'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.520'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.98'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.311'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.091.08'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '158.40'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.099'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.585'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '122.13'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9737'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09546'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.595'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.330'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06101'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02257'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.230'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.241'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6018'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.22'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.269'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5712'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.17'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.945'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.388'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '115.70'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06884'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000302'

# Set Volume(1h) column (E) cells, preserving the padded spacing
$ws.Range("E2").Value = '  +1.15%  '
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.24%  '
$ws.Range("E5").Value = '  +2.67%  '
$ws.Range("E7").Value = '  +3.99%  '
$ws.Range("E8").Value = '  +2.06%  '
$ws.Range("E9").Value = '  -1.28%  '
$ws.Range("E10").Value = '  +1.43%  '
$ws.Range("E11").Value = '  +0.94%  '
$ws.Range("E12").Value = '  +2.85%  '
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("E14").Value = '  +0.99%  '
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("E17").Value = '  +1.51%  '
$ws.Range("E18").Value = '  +1.48%  '
$ws.Range("E19").Value = '  +2.39%  '
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("E21").Value = '  +0.09%  '
$ws.Range("E22").Value = '  +1.21%  '
$ws.Range("E23").Value = '  -0.17%  '
$ws.Range("E24").Value = '  +1.25%  '
$ws.Range("E25").Value = '  +1.48%  '
$ws.Range("E26").Value = '  +0.27%  '
$ws.Range("E27").Value = '  +3.06%  '
$ws.Range("E28").Value = '  +2.35%  '
$ws.Range("E29").Value = '  +1.31%  '
$ws.Range("E30").Value = '  +2.73%  '
$ws.Range("E31").Value = '  +1.20%  '
$ws.Range("E32").Value = '  +4.04%  '
$ws.Range("E33").Value = '  +2.68%  '
$ws.Range("E34").Value = '  -2.42%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("E36").Value = '  +1.41%  '
$ws.Range("E37").Value = '  +1.87%  '
$ws.Range("E38").Value = '  +0.88%  '
$ws.Range("E39").Value = '  +0.88%  '
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("E41").Value = '  +1.84%  '
$ws.Range("E42").Value = '  +0.91%  '
$ws.Range("E43").Value = '  +1.06%  '
$ws.Range("E44").Value = '  -0.51%  '
$ws.Range("E45").Value = '  +1.66%  '
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("E47").Value = '  +1.20%  '
$ws.Range("E48").Value = '  +0.70%  '
$ws.Range("E49").Value = '  +6.98%  '
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("E51").Value = '  +15.41%  '

# Clear the temporary text-number-format so no stray style remains on the Price column
$ws.Range("D2:D51").ClearFormats()
